$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Rename sheet "Лист1" -> "Nädal 3"
$ws.Name = "Nädal 3"

# Row 9 updates
$ws.Range("D9").Value = 0.029861111111111113
$ws.Range("F9").Value = 103
$ws.Range("H9").Value = "punktid 1 - 2, alustasin p. 3."
$ws.Range("J9").Value = "x"

# Row 10 updates
$ws.Range("B10").Value = 43877
$ws.Range("C10").Value = 0.55208333333333337
$ws.Range("G10").Value = "Kodutöö 3"

# Column B width / selection
$ws.Columns.Item(2).ColumnWidth = 9.25
$ws.Range("E1").Select()
